# daily auto push: 2026-02-05 14:10 UTC
# Two new readings were recorded for 2026/02/05 (Thursday). Insert them as
# two new rows right after the existing 2026/02/05 rows (row 759), pushing
# every following row down by two (old row 760 -> new row 762, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the old row 760 (shifts 760.. down to 762..)
$ws.Rows.Item(760).Insert()
$ws.Rows.Item(760).Insert()

# Fill the two newly-inserted rows. The leading apostrophe forces column A
# to be stored as literal text ("2026/02/05") instead of being auto-parsed
# into a date serial, matching the rest of the date column.
$ws.Range("A760").Value = "'2026/02/05"
$ws.Range("B760").Value = "木"
$ws.Range("C760").Value = 18
$ws.Range("D760").Value = 47

$ws.Range("A761").Value = "'2026/02/05"
$ws.Range("B761").Value = "木"
$ws.Range("C761").Value = 22
$ws.Range("D761").Value = 53
